# Comienzo de frontend fin de sesion 27-1-2026
# - Removes the duplicated "EMPRESA 2/3", "PUESTO 2/3", "PERIODO 2/3" block
#   (columns U:Z) entirely from the sheet.
# - Clears the now-stale "Puesto actual" / "EMPRESA 1" / "PERIODO 1" /
#   "PUESTO 1" values (columns K, N, O, P) for the data rows that had them,
#   leaving the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns U:Z (EMPRESA 2/PUESTO 2/PERIODO 2/EMPRESA 3/PUESTO 3/PERIODO 3)
$ws.Range("U1:Z30").EntireColumn.Delete()

# Clear stale per-person job/company/period data in columns K, N, O, P
$rows = @(2, 4, 5, 6, 7)
foreach ($r in $rows) {
    $ws.Range("K$r").Value = ""
    $ws.Range("N$r").Value = ""
    $ws.Range("O$r").Value = ""
    $ws.Range("P$r").Value = ""
}
